{"js": "// Remove the duplicated screenshot paragraph that immediately follows the\n// \"Presentaci\u00f3n del proyecto\" (Heading 2) heading. The same picture is\n// already shown under \"Introducci\u00f3n\"; this second, redundant copy (an\n// image-only paragraph) is deleted outright, as in the source diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text,style\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if ((p.style || \"\").indexOf(\"Heading 2\") !== -1 &&\n      (p.text || \"\").trim() === \"Presentaci\u00f3n del proyecto\") {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const next = target.getNext();\n  next.load(\"text\");\n  const pics = next.inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n\n  // Only remove it when it is indeed the image-only duplicate paragraph\n  // (no real text, exactly one inline picture) \u2014 the same guard a human\n  // editor would apply before deleting.\n  if ((next.text || \"\").trim() === \"\" && pics.items.length > 0) {\n    next.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the duplicated screenshot paragraph that immediately follows the\n# \"Presentaci\u00f3n del proyecto\" (Heading 2) heading. The same picture is\n# already shown under \"Introducci\u00f3n\"; this second, redundant copy (an\n# image-only paragraph) is deleted outright, as in the source diff.\n$d = $word.ActiveDocument\n\n$heading = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text -like \"*Presentaci?n del proyecto*\") {\n        $heading = $p\n        break\n    }\n}\n\nif ($heading -ne $null) {\n    $target = $heading.Next()\n    # Only remove it when it is indeed the image-only duplicate paragraph\n    # (no real text other than the paragraph mark, exactly one inline\n    # picture) \u2014 the same guard a human editor would apply before deleting.\n    if ($target -ne $null -and $target.Range.InlineShapes.Count -eq 1 -and $target.Range.Text.Trim() -eq \"\") {\n        $target.Range.Delete()\n    }\n}\n"}
